$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The match rows below got re-ordered/re-matched against different fixtures.
# Column A (the running "id" index) stays put for each physical row; every
# other column (B..AD, i.e. B through AD) moves together with the fixture it
# describes. We read the "B:AD" row-ranges first, then write them back out in
# their new order so that reads never see already-overwritten data.

function Get-RowData($rowNum) {
    return $ws.Range("B" + $rowNum + ":AD" + $rowNum).Value2
}

function Set-RowData($rowNum, $data) {
    $ws.Range("B" + $rowNum + ":AD" + $rowNum).Value2 = $data
}

# Group 1: rows 38 and 39 simply swap places.
$r38 = Get-RowData 38
$r39 = Get-RowData 39
Set-RowData 38 $r39
Set-RowData 39 $r38

# Group 2: rows 200 and 201 simply swap places.
$r200 = Get-RowData 200
$r201 = Get-RowData 201
Set-RowData 200 $r201
Set-RowData 201 $r200

# Group 3: rows 237 and 238 simply swap places.
$r237 = Get-RowData 237
$r238 = Get-RowData 238
Set-RowData 237 $r238
Set-RowData 238 $r237

# Group 4: rows 268, 269, 270, 271 are cyclically re-matched:
#   new 268 <- old 271
#   new 269 <- old 270
#   new 270 <- old 268
#   new 271 <- old 269
$r268 = Get-RowData 268
$r269 = Get-RowData 269
$r270 = Get-RowData 270
$r271 = Get-RowData 271
Set-RowData 268 $r271
Set-RowData 269 $r270
Set-RowData 270 $r268
Set-RowData 271 $r269
